$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 253.5
$ws.Range("I4").Value = 292.7
$ws.Range("J4").Value = 57.5
$ws.Range("K4").Value = 292.7
$ws.Range("L4").Value = 57.5
$ws.Range("M4").Value = -178.7
$ws.Range("N4").Value = -285.5

# Row 20
$ws.Range("H20").Value = 11976.2
$ws.Range("J20").Value = 10000
$ws.Range("L20").Value = 10000
$ws.Range("N20").Value = -10460

# Row 35
$ws.Range("H35").Value = 11976.2
$ws.Range("J35").Value = 10000
$ws.Range("L35").Value = 10000
$ws.Range("N35").Value = -10758

# Row 64
$ws.Range("H64").Value = 3446.5789
$ws.Range("J64").Value = 4048.1875
$ws.Range("L64").Value = 4048.1875
$ws.Range("N64").Value = -4544.1875

# Row 67
$ws.Range("H67").Value = 3446.5789
$ws.Range("J67").Value = 4048.1875
$ws.Range("L67").Value = 4048.1875
$ws.Range("N67").Value = -5764.1875

# Row 98
$ws.Range("H98").Value = 652.4545000000001
$ws.Range("I98").Value = 506.89474
$ws.Range("K98").Value = 506.89474
$ws.Range("M98").Value = 991.10526

# Row 122
$ws.Range("H122").Value = 652.4545000000001
$ws.Range("I122").Value = 506.89474
$ws.Range("K122").Value = 1520.68422
$ws.Range("M122").Value = 929.3157799999999

# Row 125
$ws.Range("H125").Value = 10170
$ws.Range("I125").Value = 9505
$ws.Range("J125").Value = 11500
$ws.Range("K125").Value = 85545
$ws.Range("L125").Value = 103500
$ws.Range("M125").Value = -83085
$ws.Range("N125").Value = -108420

$ws = $wb.Worksheets.Item("ARM")
# Row 122
$ws.Range("H122").Value = 2415.375
$ws.Range("I122").Value = 2322.8462
$ws.Range("K122").Value = 6968.5386
$ws.Range("M122").Value = -4518.5386

# Row 123
$ws.Range("H123").Value = 50000
$ws.Range("J123").Value = 50000
$ws.Range("L123").Value = 50000
$ws.Range("N123").Value = -59800

# Row 128
$ws.Range("H128").Value = 72900
$ws.Range("J128").Value = 72900
$ws.Range("L128").Value = 72900
$ws.Range("N128").Value = -82860

# Row 132
$ws.Range("H132").Value = 1401.875
$ws.Range("I132").Value = 1257.909
$ws.Range("K132").Value = 3773.727
$ws.Range("M132").Value = -1243.727

# Row 138
$ws.Range("H138").Value = 66454.75
$ws.Range("J138").Value = 63809.668
$ws.Range("L138").Value = 63809.668
$ws.Range("N138").Value = -74089.66800000001

$ws = $wb.Worksheets.Item("BSM")
# Row 8
$ws.Range("H8").Value = 2600
$ws.Range("I8").Value = 200
$ws.Range("J8").Value = 5000
$ws.Range("K8").Value = 200
$ws.Range("L8").Value = 5000
$ws.Range("M8").Value = -60
$ws.Range("N8").Value = -5280

# Row 62
$ws.Range("H62").Value = 25000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""

# Row 65
$ws.Range("H65").Value = 25000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""

# Row 105
$ws.Range("H105").Value = 80453.16
$ws.Range("J105").Value = 4561.5
$ws.Range("L105").Value = 4561.5
$ws.Range("N105").Value = -8055.5

# Row 107
$ws.Range("H107").Value = 1650.8422
$ws.Range("I107").Value = 1486.625
$ws.Range("K107").Value = 1486.625
$ws.Range("M107").Value = 433.375

# Row 134
$ws.Range("H134").Value = 10449.571
$ws.Range("I134").Value = 10524.5
$ws.Range("K134").Value = 31573.5
$ws.Range("M134").Value = -29038.5

$ws = $wb.Worksheets.Item("CRP")
# Row 59
$ws.Range("H59").Value = 86598.664
$ws.Range("J59").Value = 122398
$ws.Range("L59").Value = 122398
$ws.Range("N59").Value = -124688

# Row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = ""

# Row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = ""

$ws = $wb.Worksheets.Item("CUL")
# Row 10
$ws.Range("H10").Value = 268.33334
$ws.Range("I10").Value = 201.81818
$ws.Range("J10").Value = 1000
$ws.Range("K10").Value = 605.4545400000001
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = -466.4545400000001
$ws.Range("N10").Value = -3278

$ws = $wb.Worksheets.Item("GSM")
# Row 93
$ws.Range("H93").Value = 17787
$ws.Range("J93").Value = 17787
$ws.Range("L93").Value = 17787
$ws.Range("N93").Value = -21531

# Row 122
$ws.Range("H122").Value = 2516253.5
$ws.Range("I122").Value = 3349668.8
$ws.Range("K122").Value = 10049006.4
$ws.Range("M122").Value = -10046556.4

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 3679.1365
$ws.Range("I46").Value = 2740
$ws.Range("J46").Value = 3887.8333
$ws.Range("K46").Value = 2740
$ws.Range("L46").Value = 3887.8333
$ws.Range("N46").Value = -4263.8333
$ws.Range("M46").Value = -2552

# Row 51
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").Value = ""

# Row 68
$ws.Range("H68").Value = 5901.5
$ws.Range("I68").Value = 5901.5
$ws.Range("K68").Value = 5901.5
$ws.Range("M68").Value = -5152.5

# Row 71
$ws.Range("H71").Value = 5901.5
$ws.Range("I71").Value = 5901.5
$ws.Range("K71").Value = 29507.5
$ws.Range("M71").Value = -25763.5

# Row 122
$ws.Range("H122").Value = 54549292
$ws.Range("I122").Value = 90912750
$ws.Range("K122").Value = 272738250
$ws.Range("M122").Value = -272735800

# Row 136
$ws.Range("H136").Value = 2552.6296
$ws.Range("I136").Value = 2982.647
$ws.Range("J136").Value = 1821.6
$ws.Range("K136").Value = 8947.940999999999
$ws.Range("L136").Value = 5464.799999999999
$ws.Range("M136").Value = -6397.940999999999
$ws.Range("N136").Value = -10564.8

$ws = $wb.Worksheets.Item("WVR")
# Row 70
$ws.Range("H70").Value = 53500
$ws.Range("J70").Value = 53500
$ws.Range("L70").Value = 53500
$ws.Range("N70").Value = -54130

# Row 73
$ws.Range("H73").Value = 53500
$ws.Range("J73").Value = 53500
$ws.Range("L73").Value = 53500
$ws.Range("N73").Value = -55684

# Row 86
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").Value = ""

# Row 89
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").Value = ""

# Row 107
$ws.Range("H107").Value = 1554.4615
$ws.Range("I107").Value = 699.8333
$ws.Range("J107").Value = 2287
$ws.Range("K107").Value = 2099.4999
$ws.Range("L107").Value = 6861
$ws.Range("M107").Value = -179.4998999999998
$ws.Range("N107").Value = -10701

# Row 122
$ws.Range("H122").Value = 2534.7
$ws.Range("I122").Value = 2091.1667
$ws.Range("J122").Value = 3200
$ws.Range("K122").Value = 6273.500100000001
$ws.Range("L122").Value = 9600
$ws.Range("M122").Value = -3823.500100000001
$ws.Range("N122").Value = -14500

# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = ""
